$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Sensitive to assumed inflation process?" row (row 2):
# the DE and DENI columns (D2, E2) flip from "No" to "Yes"
$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = "Yes"

# Update the "Sensitive to two-step or joint estimate?" row (row 4):
# the DE column (D4) flips from "No" to "Yes"
$ws.Range("D4").Value = "Yes"
